$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes existing rows 10-38 down to 11-39)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44525
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112031
$ws.Range("G10").Value = "Poroto verde"
$ws.Range("H10").Value = "Magnum"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 29000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 29500
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 1180
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
